$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the "tanggalLahir" header (was "tangalLahir").
# Re-assigning the text causes the shared-string table to drop the old
# unused entry and append the corrected string at the end of the table.
$ws.Range("I1").Value = "tanggalLahir"

# Move the active cell / selection from F5 to F9.
$ws.Range("F9").Select()
